# Apply "latest DAP file" updates to the total_bcq_nomination sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (HOUR 1)
$ws.Range("B2").Value = 12500
$ws.Range("C2").Value = 10000
$ws.Range("E2").Value = 0
$ws.Range("F2").Value = 42500

# Row 3 (HOUR 2)
$ws.Range("D3").Value = 0
$ws.Range("F3").Value = 22500

# Row 6 (HOUR 5)
$ws.Range("D6").Value = 0
$ws.Range("F6").Value = 22500

# Row 7 (HOUR 6)
$ws.Range("D7").Value = 0
$ws.Range("F7").Value = 22500

# Row 8 (HOUR 7)
$ws.Range("D8").Value = 0
$ws.Range("F8").Value = 22500

# Row 9 (HOUR 8)
$ws.Range("D9").Value = 0
$ws.Range("F9").Value = 22500

# Row 23 (HOUR 22)
$ws.Range("B23").Value = 24000
$ws.Range("F23").Value = 74000

# Row 24 (HOUR 23)
$ws.Range("B24").Value = 21000
$ws.Range("F24").Value = 71000

# Row 25 (HOUR 24)
$ws.Range("B25").Value = 17500
$ws.Range("F25").Value = 67500
